$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the meeting length for the 3/15 entry (row 4) to include the end time
$ws.Range("B4").Value = "10:50am - 11:15"

# Fill in the meeting overview for the 3/15 entry (row 4); match formatting of
# the neighboring "Absent" cell (plain text style) before setting the value
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = "CRM, Sprint plan, High Level Design, Project requirements docs"

# Add a new row 5 for the 3/29 meeting, copying formats from row 4 first so the
# new cells line up with the existing style (date format on A, plain text on
# the rest) instead of creating brand-new styles
$ws.Range("A4:E4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("A5").Value = (Get-Date -Year 2018 -Month 3 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B5").Value = "11:15-11:30"
$ws.Range("C5").Value = "Nicole, Feiyu, Nicole"
$ws.Range("D5").Value = "Jacob"
$ws.Range("E5").Value = "Unit testing, System Testing, Bug list, Performance document "
